$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new columns -------------------------------------------------
# Insert 3 columns before the current "Negative Regulators" column (G)
# to hold the new "Positive *" fields.
$ws.Range("G1:I1").EntireColumn.Insert()

# Insert 3 more columns before the current "Value Type" column (now K)
# to hold the new "Negative *" fields.
$ws.Range("K1:M1").EntireColumn.Insert()

# --- Fill in the new header values --------------------------------------
$ws.Range("G1").Value = "Positive Connection Type"
$ws.Range("H1").Value = "Positive Mechanism"
$ws.Range("I1").Value = "Positive Site"

$ws.Range("K1").Value = "Negative Connection Type"
$ws.Range("L1").Value = "Negative Mechanism"
$ws.Range("M1").Value = "Negative Site"

# --- Match formatting of the new cells to their neighbours --------------
# G1:I1 should look like F1 (wrapped header style)
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:I1").PasteSpecial(-4122) | Out-Null

# K1:M1 should look like the rest of the plain headers (e.g. J1)
$ws.Range("J1").Copy() | Out-Null
$ws.Range("K1:M1").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Columns F:I share the wrap-text column style used previously only by F
$ws.Range("F1:I1").EntireColumn.WrapText = $true

# --- View / selection state ----------------------------------------------
$ws.Range("K1:M1").Select()
$ws.Application.ActiveWindow.ScrollColumn = $ws.Range("L1").Column
